$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells (D/E columns) keep their literal formatting (no numeric
# auto-conversion / trailing-zero stripping) by forcing Text number format
# before assigning string values, matching the original inlineStr cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.571.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.457.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.96"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.905.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.570.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.463.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.17"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.23%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +15.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "636.60"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0968"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -20.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.17%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.06%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.367"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "150.41"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.78%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -23.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.56"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.605"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0906"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.79%  "
